# "dodano skrypty na przywrócenie identity na ID"
#
# The "CURRENTLY" section (a blank separator row + its 5 items, one of
# which duplicated text already present under TO_DO) is removed. Two of
# the remaining TO_DO-only items are dropped as well, and the two
# genuinely-new items that used to live under CURRENTLY ("Potestować..."
# and "Potworzyć...") are kept, sliding up to the bottom of the now
# shorter TO_DO list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite the TO_DO list (row 1 keeps its header formatting/value).
$ws.Range("A2").Value = "podzielić na pliki "
$ws.Range("A3").Value = "Swagger"
$ws.Range("A4").Value = "Inżynierka"
$ws.Range("A5").Value = "Potestować wszystkie endpointy na wszelkie sposoby"
$ws.Range("A6").Value = "Potworzyć gotowe zbiory call’i w PostManie"

# Drop the old CURRENTLY block entirely (separator row 7 through row 12).
$ws.Rows("7:12").Delete()

# Match the author's final selection.
$ws.Range("F3:F4").Select()
